$d = $word.ActiveDocument
$newText = "Perioadele campaniei din Perseus: 16-25 ianuarie, 7-16 noiembrie, 6-15 decembrie"
$marker = "Perioadele campaniei"

# Collect indices first (mutating while iterating the live collection is unsafe).
$targets = New-Object System.Collections.ArrayList
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*$marker*") {
        [void]$targets.Add($i)
    }
}

foreach ($i in $targets) {
    $p = $d.Paragraphs($i)
    $full = $p.Range
    # Range covering the paragraph's content only (exclude the trailing
    # paragraph-mark character) so the paragraph mark / pPr (and any
    # sectPr riding on it) survive untouched.
    $contentRng = $d.Range($full.Start, $full.End - 1)
    $contentRng.Delete()
    # Insert fresh, completely unformatted text (no inherited rPr) at the
    # now-empty paragraph start.
    $insertPoint = $d.Range($full.Start, $full.Start)
    $insertPoint.InsertAfter($newText)
}
